# Auto-generated: apply cryptos list price/volume updates (2024-06-30 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.530.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.387.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.966.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.382.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.580.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000112"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  +7.98%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.422.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.469.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.01%  "
$ws.Range("E51").Value = "  -2.32%  "
